# AstroLib Master Layout: rename ".Observing" -> ".Observer" and rework its content
# (Doxygen adds; start on AstroLib.Observer)

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Work on the ".Observing" sheet (to become ".Observer")
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item(".Observing")

# Remove the old body rows (5 through 13) completely (content + formatting),
# so no stray row-height / style residue is left behind.
$ws.Range("A5:D13").EntireRow.Delete()

# NOTE: the brand-new shared strings introduced by this edit must be written
# in the same order the original author typed them, so that the shared
# string table ends up in the same order. That order is:
#   Earth site definition: relatively fixed
#   Earth site's atmosphere; may be variable
#   project: AstroLib.Observer
#   namespace: AstroLib.Observer
#   Shelter

$ws.Range("D5").Value = "Earth site definition: relatively fixed"
$ws.Range("D6").Value = "Earth site's atmosphere; may be variable"

# --- Row 1: page title -------------------------------------------------
$ws.Range("A1").Value = "project: AstroLib.Observer"
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Font.Size = 16

# --- Row 4: namespace header --------------------------------------------
$ws.Range("A4").Value = "namespace: AstroLib.Observer"
$ws.Range("A4").Font.Bold = $true
$ws.Range("A4").Font.Size = 14

# --- Row 9: Shelter name (new string) -------------------------------------
$ws.Range("C9").Value = "Shelter"
$ws.Range("C9").Font.Bold = $true
$ws.Range("C9").HorizontalAlignment = -4131   # xlLeft

# --- Row 2: descriptive text -------------------------------------------
$ws.Range("A2").Value = "Defined by user; each includes numerous timings for planning."
$ws.Range("A2").Font.Italic = $true
$ws.Range("A2").HorizontalAlignment = -4131   # xlLeft

# --- Row 5: Site ----------------------------------------------------------
$ws.Range("B5").Value = "classes:"
$ws.Range("B5").Font.Italic = $true
$ws.Range("B5").HorizontalAlignment = -4152   # xlRight

$ws.Range("C5").Value = "Site"
$ws.Range("C5").Font.Bold = $true
$ws.Range("C5").HorizontalAlignment = -4131   # xlLeft

# --- Row 6: Atmosphere ------------------------------------------------
$ws.Range("C6").Value = "Atmosphere"
$ws.Range("C6").Font.Bold = $true
$ws.Range("C6").HorizontalAlignment = -4131   # xlLeft

# --- Row 7: Telescope ---------------------------------------------------
$ws.Range("B7").Font.Italic = $true
$ws.Range("B7").HorizontalAlignment = -4152   # xlRight

$ws.Range("C7").Value = "Telescope"
$ws.Range("C7").Font.Bold = $true
$ws.Range("C7").HorizontalAlignment = -4131   # xlLeft

$ws.Range("D7").Value = "includes focal length, aperture, range of motion"

# --- Row 8: Camera -------------------------------------------------------
$ws.Range("B8").Font.Italic = $true
$ws.Range("B8").HorizontalAlignment = -4152   # xlRight

$ws.Range("C8").Value = "Camera"
$ws.Range("C8").Font.Bold = $true
$ws.Range("C8").HorizontalAlignment = -4131   # xlLeft

$ws.Range("D8").Value = "includes filters, FOV, pixel scale, etc"

# --- Row 9 (cont'd): Shelter row ------------------------------------------
$ws.Range("B9").HorizontalAlignment = -4131   # xlLeft

$ws.Range("D9").Value = "dome or roof"

# --- Rename the sheet ----------------------------------------------------
$ws.Name = ".Observer"

# --- Selection / activation ----------------------------------------------
$ws.Range("A1").Select()
$ws.Activate()
$ws.Range("C10").Select()
